$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Rows 3-26 each need:
#  - column H (PERIOD TO EXPIRE) decremented by 1
#  - column I (LAST UPDATE) changed from 03-Nov-2025 to 04-Nov-2025
$firstRow = 3
$lastRow = 26

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $hCell = $ws.Cells.Item($row, 8)
    $hCell.Value2 = $hCell.Value2 - 1
}

# Setting a plain date-like string via .Value/.Formula causes the engine to
# auto-parse it as a date serial (and assign a new date number-format style),
# which would not match the original inline-string / style. Instead, build
# a text formula for each cell, then convert the range to static values via
# Copy + PasteSpecial(xlPasteValues) so the result is a literal text value
# (not a formula, not a parsed date) while the cell's existing style is left
# untouched.
for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 9).Formula = '="04-Nov-2025"'
}
$dateRange = $ws.Range("I$firstRow`:I$lastRow")
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)
$excel.CutCopyMode = 0
